# Update "相談件数" (consultation counts) sheet with one additional day of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Insert a new row above the current last data row (117), pushing the
# trailing "note" row down to 118. Insert with xlShiftDown-like behaviour
# using EntireRow.Insert(), which also shifts formatting/formulas down.
$ws.Rows.Item(117).Insert()

# Fill the newly inserted row 117 with the new day's figures. The inserted
# row already inherited the correct number formats, so just set values.
$ws.Cells.Item(117, 1).Value = 43972
$ws.Cells.Item(117, 2).Value = 119
$ws.Cells.Item(117, 3).Value = 38465
$ws.Cells.Item(117, 4).Value = 36
$ws.Cells.Item(117, 5).Value = 7764

# Keep the sheet selection consistent with the shifted note row.
$ws.Range("B118").Select()

# Extend the workbook-level Print_Area defined name to cover the new row.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area") {
        $n.RefersTo = "=相談件数!`$A`$1:`$E`$118"
    }
}
